# B2.2 workbook update
# - fills in the "Widerstand" sheet (sheet3) with new r/b/l/d measurement
#   data and the derived Bragg-additivity / lattice-parameter calculations
# - touches up sheet view state (active sheet, scroll position, selection)
#   to match the author's last-saved session

$wb = $excel.ActiveWorkbook

$wsTabelle1   = $wb.Worksheets.Item("Tabelle1")
$wsXRD        = $wb.Worksheets.Item("XRD")
$wsWiderstand = $wb.Worksheets.Item("Widerstand")

# ---------------------------------------------------------------------------
# Widerstand sheet: headers (row 2)
# ---------------------------------------------------------------------------
$wsWiderstand.Range("G2").Value = "r"
$wsWiderstand.Range("H2").Value = "dr"
$wsWiderstand.Range("I2").Value = "b"
$wsWiderstand.Range("J2").Value = "db"
$wsWiderstand.Range("K2").Value = "l"
$wsWiderstand.Range("L2").Value = "dl"
$wsWiderstand.Range("M2").Value = "d"
$wsWiderstand.Range("N2").Value = "dd"
$wsWiderstand.Range("P2").Value = "dp"
$wsWiderstand.Range("Q2").Value = "p"

# ---------------------------------------------------------------------------
# Widerstand sheet: measurement rows 3-5 (N wires = 2,3,4)
# ---------------------------------------------------------------------------
$wsWiderstand.Range("F3").Value = 2
$wsWiderstand.Range("G3").Value = 3.8999999999999999E-4
$wsWiderstand.Range("H3").Value = 9.9000000000000005E-7
$wsWiderstand.Range("I3").Value = 5.1000000000000004E-3
$wsWiderstand.Range("J3").Value = 3.0000000000000001E-5
$wsWiderstand.Range("K3").Value = 6.7000000000000002E-3
$wsWiderstand.Range("L3").Value = 5.0000000000000002E-5
$wsWiderstand.Range("M3").Value = 2.0000000000000001E-4
$wsWiderstand.Range("N3").Value = 5.0000000000000002E-5
$wsWiderstand.Range("P3").Formula = "=SQRT((J3*M3*G3/K3)^2+(I3*G3*N3/K3)^2+(I3*M3*H3/K3)^2+(I3*M3*G3*L3/K3^2)^2)"
$wsWiderstand.Range("Q3").Formula = "=I3*M3*G3/K3"

$wsWiderstand.Range("F4").Value = 3
$wsWiderstand.Range("G4").Value = 2.5999999999999998E-4
$wsWiderstand.Range("H4").Value = 1.5E-6
$wsWiderstand.Range("I4").Value = 5.1999999999999998E-3
$wsWiderstand.Range("J4").Value = 3.0000000000000001E-5
$wsWiderstand.Range("K4").Value = 4.3E-3
$wsWiderstand.Range("L4").Value = 5.0000000000000002E-5
$wsWiderstand.Range("M4").Value = 2.0000000000000001E-4
$wsWiderstand.Range("N4").Value = 5.0000000000000002E-5
$wsWiderstand.Range("P4").Formula = "=SQRT((J4*M4*G4/K4)^2+(I4*G4*N4/K4)^2+(I4*M4*H4/K4)^2+(I4*M4*G4*L4/K4^2)^2)"
$wsWiderstand.Range("Q4").Formula = "=I4*M4*G4/K4"

$wsWiderstand.Range("F5").Value = 4
$wsWiderstand.Range("G5").Value = 2.5999999999999999E-3
$wsWiderstand.Range("H5").Value = 9.0999999999999997E-7
$wsWiderstand.Range("I5").Value = 5.1000000000000004E-3
$wsWiderstand.Range("J5").Value = 3.0000000000000001E-5
$wsWiderstand.Range("K5").Value = 4.1000000000000003E-3
$wsWiderstand.Range("L5").Value = 5.0000000000000002E-5
$wsWiderstand.Range("M5").Value = 2.0000000000000001E-4
$wsWiderstand.Range("N5").Value = 5.0000000000000002E-5
$wsWiderstand.Range("P5").Formula = "=SQRT((J5*M5*G5/K5)^2+(I5*G5*N5/K5)^2+(I5*M5*H5/K5)^2+(I5*M5*G5*L5/K5^2)^2)"
$wsWiderstand.Range("Q5").Formula = "=I5*M5*G5/K5"

# ---------------------------------------------------------------------------
# Widerstand sheet: derived results (rows 10-14)
# ---------------------------------------------------------------------------
$wsWiderstand.Range("F10").Value = "A"
$wsWiderstand.Range("G10").Formula = "=(P5-P3)/(0.25*0.75)"
$wsWiderstand.Range("I10").Value = "S"
$wsWiderstand.Range("J10").Formula = "=SQRT((Q5-Q4)/(Q5-Q3))"

$wsWiderstand.Range("F11").Value = "dA"
$wsWiderstand.Range("G11").Formula = "=1/(0.25*0.75) * SQRT(P5^2+P3^2)"
$wsWiderstand.Range("I11").Value = "dS"
$wsWiderstand.Range("J11").Formula = "=SQRT(J12+J13+J14)"

$wsWiderstand.Range("I12").Value = "a"
$wsWiderstand.Range("J12").Formula = "=(1/(2*J10*(Q5-Q3)) * P4)^2"

$wsWiderstand.Range("I13").Value = "b"
$wsWiderstand.Range("J13").Formula = "=((Q5-Q4)*P3/(2*J10*(Q5-Q3)^2))^2"

$wsWiderstand.Range("I14").Value = "x"
$wsWiderstand.Range("J14").Formula = "=(P5*(Q4-Q3)/(J10*2*(Q5-Q3)^2))^2"

# ---------------------------------------------------------------------------
# Widerstand sheet: number formatting (scientific, matches built-in 0.00E+00)
# ---------------------------------------------------------------------------
$sciRanges = "G3:N3,G4:N4,G5:N5,Q3,Q4,Q5,J12,J13,J14"
$wsWiderstand.Range($sciRanges).NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Widerstand sheet: column widths
# ---------------------------------------------------------------------------
$wsWiderstand.Columns.Item(7).ColumnWidth = 11.1667
$wsWiderstand.Columns.Item(16).ColumnWidth = 11.1667

# ---------------------------------------------------------------------------
# Sheet view / selection state
# ---------------------------------------------------------------------------
$wsTabelle1.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsTabelle1.Range("I12").Select()

$wsXRD.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 2
$wsXRD.Range("O55").Select()

$wsWiderstand.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$wsWiderstand.Range("G12").Select()
